# feat: update version to 1.1.4 and enhance template download functionality
#
# The user-template.xlsx ships with a sample row (Name: "Zemu",
# Email: "zemu@gmail.com", mailto hyperlink) that was only there for
# testing. Clean the template so only the header row (Email/Name) is
# left; the second row stays in the sheet (still carrying the
# Hyperlink cell style on B2) but is otherwise empty, ready to be
# filled in by the "download template" feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the sample values from row 2.
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# The mailto: hyperlink was only attached to the sample e-mail address;
# remove it now that the cell is blank.
$ws.Hyperlinks.Delete()

# Leave the selection where the author left it when saving.
$ws.Range("B2").Select()
